$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "97.236.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +5.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.115.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.23%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.95%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "610.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.40%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.23%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.384"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.110.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.788"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.40%  "
$ws.Range("E12").Value = "  -0.48%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "96.772.28"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000241"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "33.92"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.05%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.66%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.686.40"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.23%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.103.45"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.77%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.58"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -6.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "517.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +17.97%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.51"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.34%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.35%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.0000195"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.87%  "
$ws.Range("E25").Value = "  -1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "86.80"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.39%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.63"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "3.279.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.14%  "
$ws.Range("E29").Value = "  +0.23%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.239"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.57%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.174"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.53%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.125"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.50%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.997"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "9.02"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.88%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "26.69"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.48%  "
$ws.Range("E36").Value = "  -5.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "7.33"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -9.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "485.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +4.16%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.54%  "
$ws.Range("E40").Value = "  +1.53%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.437"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.51%  "
$ws.Range("E42").Value = "  -2.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.64"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -9.01%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.21"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.61%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "161.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.03%  "
$ws.Range("E47").Value = "  +3.70%  "
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.65%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.38"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.64%  "
$ws.Range("B51").Value = "VeChain"
$ws.Range("C51").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0322"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.48%  "
